# Commit: "Hab mal ne Anleitung für VPN und das Training geschrieben,
#          zu finden unter meetings->getting started"
#
# The only content change on the (single) slide is the sub-title text,
# which is shortened from "Weekly Topic:" to "Topic:".

$p = $ppt.ActivePresentation
$s = $p.Slides.Item(1)

# Locate the sub-title placeholder robustly (rather than hard-coding an
# index) by scanning for the placeholder that currently holds the
# "Weekly Topic:" text.
for ($i = 1; $i -le $s.Shapes.Count; $i++) {
    $shape = $s.Shapes.Item($i)
    if ($shape.HasTextFrame) {
        $tr = $shape.TextFrame.TextRange
        if ($tr.Text -eq "Weekly Topic:") {
            $tr.Text = "Topic:"
        }
    }
}
